# Weekly update: insert a new price record for "Arveja Verde"
# (Vega Modelo de Temuco) above the existing row 31, pushing the
# historical rows (old 31..89) down to (32..90).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at row 31; Excel shifts rows 31:89 down to 32:90
# and the sheet's used range grows from A1:R89 to A1:R90 automatically.
$ws.Rows.Item(31).Insert()

# Populate the newly inserted row 31 with this week's record.
$ws.Cells.Item(31, 1).Value = 10
$ws.Cells.Item(31, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(31, 3).Value = "La Araucanía"
$ws.Cells.Item(31, 4).Value = 44581
$ws.Cells.Item(31, 5).Value = 9
$ws.Cells.Item(31, 6).Value = 100112022
$ws.Cells.Item(31, 7).Value = "Arveja Verde"
$ws.Cells.Item(31, 8).Value = "Sin especificar"
$ws.Cells.Item(31, 9).Value = "Primera"
$ws.Cells.Item(31, 10).Value = 40
$ws.Cells.Item(31, 11).Value = 26000
$ws.Cells.Item(31, 12).Value = 26000
$ws.Cells.Item(31, 13).Value = 26000
$ws.Cells.Item(31, 14).Value = "`$/saco 25 kilos"
$ws.Cells.Item(31, 15).Value = "Región de La Araucanía"
$ws.Cells.Item(31, 16).Value = 1040
$ws.Cells.Item(31, 17).Value = 25
$ws.Cells.Item(31, 18).Value = "Hortaliza"
